$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FedEx shipment tracking numbers for column P (ShipmentTracking), rows 2-26
$trackingNumbers = @(
    "320018594180",
    "320018594190",
    "320018594227",
    "320018594249",
    "320018594282",
    "320018594308",
    "320018594330",
    "320018594352",
    "320018594385",
    "320018594400",
    "320018594444",
    "320018594466",
    "320018594499",
    "320018594514",
    "320018594547",
    "320018594569",
    "320018594606",
    "320018594628",
    "320018594650",
    "320018594672",
    "320018594709",
    "320018594710",
    "320018594720",
    "320018594731",
    "320018594742"
)

$row = 2
foreach ($tn in $trackingNumbers) {
    # Prefix with an apostrophe so Excel stores the numeric-looking value as text
    # (matches the existing shared-string/text type already used in column P).
    $ws.Range("P" + $row).Value = "'" + $tn
    $row = $row + 1
}

# Drop the quote-prefix formatting Excel applied so the cell style matches the
# rest of the un-styled data cells in this column.
$ws.Range("P2:P26").ClearFormats()
